$d = $word.ActiveDocument

# Remove the three image paragraphs (each paragraph consists solely of an
# inline drawing). We delete them from the end backwards so the indices /
# ranges of the earlier paragraphs stay valid.
for ($i = 3; $i -ge 1; $i--) {
    $para = $d.Paragraphs.Item($i)
    $para.Range.Delete()
}

# Fix up the typos in the remaining instructional text. We locate each
# target phrase with Find (leaving the match string in place, i.e. an
# empty replacement) and then overwrite the narrowed Range's .Text
# directly -- this avoids Find/Replace's "smart quote" autocorrection,
# which would otherwise turn the literal straight apostrophe into a
# curly one.
$rng1 = $d.Content
$rng1.Find.Execute("Do not open this booklet until you are told fo do.so.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng1.Text = "Do not open this bookiet until you are told fo do'so."

$rng2 = $d.Content
$rng2.Find.Execute("Follow alt instructions carefully.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng2.Text = "Foltow all instructions carefully."

$d.Save()
